$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new staff rows (33 and 34) ---------------------------------
# Values are entered in the same order the source workbook's shared-string
# table was built in (Name col, then Designation, then Photo, then unique_id)
# so that new shared strings line up with the target workbook.
$ws.Range("A33").Value = "DHANALAKSHMI G"
$ws.Range("A34").Value = "SUBRAMANIAN N"

$ws.Range("B33").Value = "Lab Instructor"
$ws.Range("B34").Value = "Lab Instructor"

$ws.Range("C34").Value = "/static/images/profile_photos/009/VEC-009-05-4.webp"
$ws.Range("C33").Value = "/static/images/profile_photos/009/VEC-009-05-5.webp"

$ws.Range("J33").Value = "VEC-009-05-5"
$ws.Range("J34").Value = "VEC-009-05-4"

# --- Remove the stale (invisible) border formatting from A26:A32 -----------
$ws.Range("A26:A32").Borders.LineStyle = -4142

# --- Apply a thin box border around the new Name/Designation cells ---------
$ws.Range("A33:A34").Borders.LineStyle = 1
$ws.Range("A33:A34").Borders.Weight = 2

$ws.Range("B33:B34").Borders.LineStyle = 1
$ws.Range("B33:B34").Borders.Weight = 2

# --- Explicit black font colour for the new Designation cells --------------
$ws.Range("B33:B34").Font.Color = 0

# --- Restore the selection/active cell to the newly added row --------------
$ws.Range("A35").Select()
